$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 4182.119796911906
$ws.Range("C3").Value = 4177.315687199358
$ws.Range("C4").Value = 4148.239638268481
$ws.Range("C5").Value = 4148.239638268481
$ws.Range("C6").Value = 4141.828649114111
$ws.Range("C7").Value = 4139.045239324526
$ws.Range("C8").Value = 4086.578178276978
$ws.Range("C9").Value = 4080.148293248155
$ws.Range("C10").Value = 4080.148293248155
$ws.Range("C11").Value = 4080.148293248155
$ws.Range("C12").Value = 4080.148293248155
